$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "IC1, IC2"

$ws.Range("B5").Select()
